# Created new test case for publons module
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# Row 31: merge in the new OPQA ticket / description for the existing alias-account test case
$ws.Range("B31").Value = "OPQA-5993||OPQA-5995"
$ws.Range("C31").Value = "Verify User able to add alias account after click on ""Add email address"" and  Verify that user received activation link,after click on activation link user is created and navigating to account setting page"
$ws.Rows.Item(31).RowHeight = 30

# Row 32: brand new test case for switching the primary account
$ws.Range("C32").Value = "Verify user can switch primary account within added alias account"
$ws.Range("B32").Value = "OPQA-5998"
$ws.Range("D32").Value = "Y"

# Update selection to reflect where the author left off editing
$ws.Range("B32").Select()
